$d = $word.ActiveDocument

# 1. "estagiário durante mais de 1 ano e até o presente momento nesta empresa
#     maravilhosa, com excelente..." -> "...mais de 1 ano nesta ótima empresa,
#     com excelente..."
$d.Content.Find.Execute(
    "durante mais de 1 ano e até o presente momento nesta empresa maravilhosa, com excelente",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "durante mais de 1 ano nesta ótima empresa, com excelente", 2) | Out-Null

# 2. Turn the manual line break right after "Google Play Console. " into a
#    real paragraph break, splitting the sentence about DevOps experience
#    from "Essas atividades facilitaram..." into its own paragraph. Only the
#    first manual line break in the document (this one) is targeted.
$d.Content.Find.Execute(
    "^l", $true, $false, $false, $false, $false, $true, 1, $false,
    "^p", 1) | Out-Null

# 3. Move the "_GoBack" bookmark: it used to sit inside "maior parte" in the
#    hobbies paragraph; it now marks the last edit point, between ".Net MVC, "
#    and "SQL, " in the paragraph we just split off above. Adding a bookmark
#    with the same name automatically relocates it (Word keeps only one).
$rng = $d.Content
$rng.Find.Execute("Net MVC, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null

# 4. Now that the bookmark no longer splits "maior p" / "arte..." apart,
#    re-write that sentence as one contiguous run of text ("maior parte").
$d.Content.Find.Execute(
    "aproveito a maior parte do meu tempo livre praticando esportes radicais, como skate e mountain",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "aproveito a maior parte do meu tempo livre praticando esportes radicais, como skate e mountain", 2) | Out-Null
